$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 2182468.8
$ws.Range("I9").Value = 150
$ws.Range("J9").Value = 2400700.5
$ws.Range("K9").Value = 150
$ws.Range("L9").Value = 2400700.5
$ws.Range("M9").Value = 19
$ws.Range("N9").Value = -2401038.5
$ws.Range("H40").Value = 2613.9473
$ws.Range("I40").Value = 2193.3333
$ws.Range("K40").Value = 2193.3333
$ws.Range("M40").Value = -2018.3333
$ws.Range("H98").Value = 2361.9048
$ws.Range("I98").Value = 1980.5
$ws.Range("K98").Value = 1980.5
$ws.Range("M98").Value = -482.5
$ws.Range("H122").Value = 2361.9048
$ws.Range("I122").Value = 1980.5
$ws.Range("K122").Value = 5941.5
$ws.Range("M122").Value = -3491.5
$ws.Range("H125").Value = 2749
$ws.Range("J125").Value = 2749
$ws.Range("L125").Value = 24741
$ws.Range("N125").Value = -29661
$ws.Range("H137").Value = 17859628
$ws.Range("J137").Value = 2799.9092
$ws.Range("L137").Value = 8399.7276
$ws.Range("N137").Value = -13499.7276
$ws.Range("H138").Value = 2964.532
$ws.Range("I138").Value = 1459.9412
$ws.Range("J138").Value = 3817.1333
$ws.Range("K138").Value = 4379.8236
$ws.Range("L138").Value = 11451.3999
$ws.Range("M138").Value = 760.1764000000003
$ws.Range("N138").Value = -21731.3999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19334.12
$ws.Range("I32").Value = 19968.299
$ws.Range("K32").Value = 19968.299
$ws.Range("M32").Value = -19681.299
$ws.Range("H61").Value = 4765.222
$ws.Range("I61").Value = 3650.5715
$ws.Range("K61").Value = 3650.5715
$ws.Range("M61").Value = -3438.5715
$ws.Range("H102").Value = 4342
$ws.Range("I102").Value = 3676.5715
$ws.Range("K102").Value = 3676.5715
$ws.Range("M102").Value = -2054.5715
$ws.Range("H136").Value = 4765.222
$ws.Range("I136").Value = 3650.5715
$ws.Range("K136").Value = 10951.7145
$ws.Range("M136").Value = -8401.7145

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 52636160
$ws.Range("I31").Value = 166667540
$ws.Range("J31").Value = 6294.4614
$ws.Range("K31").Value = 166667540
$ws.Range("L31").Value = 6294.4614
$ws.Range("M31").Value = -166667245
$ws.Range("N31").Value = -6884.4614
$ws.Range("H34").Value = 52636160
$ws.Range("I34").Value = 166667540
$ws.Range("J34").Value = 6294.4614
$ws.Range("K34").Value = 166667540
$ws.Range("L34").Value = 6294.4614
$ws.Range("M34").Value = -166667338
$ws.Range("N34").Value = -6698.4614
$ws.Range("H105").Value = 928.9091
$ws.Range("I105").Value = 921.8
$ws.Range("J105").Value = 1000
$ws.Range("K105").Value = 921.8
$ws.Range("L105").Value = 1000
$ws.Range("M105").Value = 825.2
$ws.Range("N105").Value = -4494
$ws.Range("H132").Value = 4678.3706
$ws.Range("I132").Value = 3366.1765
$ws.Range("J132").Value = 6909.1
$ws.Range("K132").Value = 10098.5295
$ws.Range("L132").Value = 20727.3
$ws.Range("M132").Value = -7568.529500000001
$ws.Range("N132").Value = -25787.3
$ws.Range("H134").Value = 7915.9414
$ws.Range("I134").Value = 7441.2144
$ws.Range("J134").Value = 10131.333
$ws.Range("K134").Value = 22323.6432
$ws.Range("L134").Value = 30393.999
$ws.Range("M134").Value = -19788.6432
$ws.Range("N134").Value = -35463.999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 1191.6666
$ws.Range("J69").Value = 1550
$ws.Range("L69").Value = 4650
$ws.Range("N69").Value = -6272
$ws.Range("H72").Value = 1191.6666
$ws.Range("J72").Value = 1550
$ws.Range("L72").Value = 13950
$ws.Range("N72").Value = -22062
$ws.Range("H92").Value = 1223.4762
$ws.Range("I92").Value = 1342.7273
$ws.Range("K92").Value = 4028.1819
$ws.Range("M92").Value = -2780.1819
$ws.Range("H131").Value = 15877537
$ws.Range("I131").Value = 47619684
$ws.Range("J131").Value = 6463.5
$ws.Range("K131").Value = 142859052
$ws.Range("L131").Value = 19390.5
$ws.Range("M131").Value = -142854012
$ws.Range("N131").Value = -29470.5
$ws.Range("H140").Value = 964.1429000000001
$ws.Range("I140").Value = 964.1429000000001
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 2892.4287
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = 2287.5713
$ws.Range("N140").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2700.2449
$ws.Range("I22").Value = 1725.7142
$ws.Range("J22").Value = 3999.6191
$ws.Range("K22").Value = 1725.7142
$ws.Range("L22").Value = 3999.6191
$ws.Range("M22").Value = -1430.7142
$ws.Range("N22").Value = -4589.6191
$ws.Range("H27").Value = 2700.2449
$ws.Range("I27").Value = 1725.7142
$ws.Range("J27").Value = 3999.6191
$ws.Range("K27").Value = 1725.7142
$ws.Range("L27").Value = 3999.6191
$ws.Range("M27").Value = -1618.7142
$ws.Range("N27").Value = -4213.6191
$ws.Range("H61").Value = 4731.72
$ws.Range("I61").Value = 5012.9565
$ws.Range("J61").Value = 1497.5
$ws.Range("K61").Value = 5012.9565
$ws.Range("L61").Value = 1497.5
$ws.Range("M61").Value = -4810.9565
$ws.Range("N61").Value = -1901.5
$ws.Range("H113").Value = 4731.72
$ws.Range("I113").Value = 5012.9565
$ws.Range("J113").Value = 1497.5
$ws.Range("K113").Value = 5012.9565
$ws.Range("L113").Value = 1497.5
$ws.Range("M113").Value = -2842.9565
$ws.Range("N113").Value = -5837.5
$ws.Range("H133").Value = 55848.668
$ws.Range("J133").Value = 55848.668
$ws.Range("L133").Value = 55848.668
$ws.Range("N133").Value = -60908.668

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 18567.5
$ws.Range("I38").Value = 30351.666
$ws.Range("K38").Value = 30351.666
$ws.Range("M38").Value = -29878.666
$ws.Range("H55").Value = 1016
$ws.Range("I55").Value = 949
$ws.Range("J55").Value = 1150
$ws.Range("K55").Value = 949
$ws.Range("L55").Value = 1150
$ws.Range("M55").Value = -672
$ws.Range("N55").Value = -1704
$ws.Range("H132").Value = 5911.3076
$ws.Range("J132").Value = 9506.5
$ws.Range("L132").Value = 28519.5
$ws.Range("N132").Value = -33579.5
$ws.Range("H136").Value = 4003.5625
$ws.Range("I136").Value = 1874.8889
$ws.Range("K136").Value = 5624.6667
$ws.Range("M136").Value = -3074.6667
